$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.911.86'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.981.07'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.57'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.12'
$ws.Range('E7').Value = '  +3.27%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +2.16%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0803'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('E11').Value = '  +0.73%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.67'
$ws.Range('E12').Value = '  +6.60%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.847'
$ws.Range('E13').Value = '  +2.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.15'
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('D15').Value = '2.269.26'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.43'
$ws.Range('E16').Value = '  +2.76%  '
$ws.Range('D17').Value = '1.969.86'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '36.787.06'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.21'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '230.44'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  +3.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.30'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '163.39'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.51'
$ws.Range('E29').Value = '  +0.42%  '
$ws.Range('E30').Value = '  +18.94%  '
$ws.Range('E31').Value = '  +2.16%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.85'
$ws.Range('E32').Value = '  +2.41%  '
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.54'
$ws.Range('E34').Value = '  +5.74%  '
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.35'
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.50'
$ws.Range('E39').Value = '  -9.98%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0973'
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.18'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0212'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.15'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '1.372.55'
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '89.92'
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '46.32'
$ws.Range('E50').Value = '  +5.83%  '
$ws.Range('D51').Value = '2.163.04'
$ws.Range('E51').Value = '  +0.98%  '
